# [external commands] - [tail(id,file)]: simulate the *NIX tail command.
#
# The "#system" sheet (first sheet in the workbook) backs several named
# ranges used for data-validation drop-downs elsewhere in the workbook.
# This commit adds a new "external" command - tail(id,file) - and a new
# "web" command - assertTextNotContains(locator,text) - to those lookup
# lists, and widens the corresponding named ranges to cover the newly
# added rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# --- "external" list (column I) --------------------------------------
# Previously I2:I4 (runJUnit / runProgram / runProgramNoWait). Append the
# new command as the next row, in alphabetical order.
$ws.Range("I5").Value = "tail(id,file)"
$wb.Names.Item("external").RefersTo = "='#system'!`$I`$2:`$I`$5"

# --- "web" list (column Y) --------------------------------------------
# Previously Y2:Y127, alphabetically sorted. The new command sorts right
# before "assertTextNotPresent(text)" (row 39), so insert a cell there -
# shifting Y39:Y127 down to Y40:Y128 - and populate the new row.
$ws.Range("Y39").Insert()
$ws.Range("Y39").Value = "assertTextNotContains(locator,text)"
$wb.Names.Item("web").RefersTo = "='#system'!`$Y`$2:`$Y`$128"
